$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 29, shifting existing rows 29-86 down to 30-87.
$ws.Rows("29").Insert()

# Populate the newly inserted row 29 with the new data record.
$ws.Cells.Item(29, 1).Value = 10
$ws.Cells.Item(29, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(29, 3).Value = "La Araucanía"
$ws.Cells.Item(29, 4).Value = 44973
$ws.Cells.Item(29, 5).Value = 9
$ws.Cells.Item(29, 6).Value = "Fruta"
$ws.Cells.Item(29, 7).Value = 100108
$ws.Cells.Item(29, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(29, 9).Value = 100108004
$ws.Cells.Item(29, 10).Value = "Papaya"
$ws.Cells.Item(29, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(29, 12).Value = "Primera"
$ws.Cells.Item(29, 13).Value = 35
$ws.Cells.Item(29, 14).Value = 42000
$ws.Cells.Item(29, 15).Value = 42000
$ws.Cells.Item(29, 16).Value = 42000
$ws.Cells.Item(29, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(29, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(29, 19).Value = 2800
$ws.Cells.Item(29, 20).Value = 15
